$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Engineered Features")

# New engineered-feature rows (GeoSpatial v3: landmark/transit features + crime data).
# "source" (column C) is intentionally entered in a second pass below, matching how
# the sheet was actually authored (new rows typed first, sources backfilled after).
$rows = @(
    @("dist_navy_pier_miles",   "float", "haversine_dist(listing, navy_pier)",    "Straight-line distance to Navy Pier in miles"),
    @("is_landmark_proximate",  "int",   "(dist_to_landmarks <= 1.0mi).astype(int)", "Binary flag for listings within 1 mile of major Chicago landmarks"),
    @("transit_north_premium",  "float", "(1 / (dist_navy_pier + 1)) * lat_norm", "Interaction between landmark proximity and North Side latitude"),
    @("is_premium_transit_hub", "float", "is_rail_accessible * lat_norm",         "Weighted score for rail access (within 0.5mi) scaled by North Side prestige"),
    @("shoreline_gravity",      "float", "1 / (dist_to_shoreline**2 + 1)",        "Exponentially decaying score based on distance to Lake Michigan"),
    @("dist_to_beach_miles",    "float", "min(dist_to_all_public_beaches)",       "Distance to the nearest public beach access point"),
    @("crime_density_total",    "float", "count(crimes) / area_sq_miles",         "Neighborhood-level crime incidents per square mile"),
    @("nb_homicides_half_mile", "int",   "count(homicides within 0.5mi buffer)",  "Localized count of homicides within walking distance")
)

$startRow = 12
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    # The crime source is typed inline with its row; the other new rows' sources
    # are filled in afterwards.
    if ($r -eq 18) {
        $ws.Cells.Item($r, 3).Value = "Crimes (Chicago Data Portal)"
    }
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = "feature"
}

# Backfill the "source" column for the remaining new rows.
$ws.Cells.Item(12, 3).Value = "listings.csv, Chicago Data Portal"
$ws.Cells.Item(13, 3).Value = "listings.csv, Chicago Data Portal"
$ws.Cells.Item(14, 3).Value = "listings.csv, Chicago Data Portal"
$ws.Cells.Item(15, 3).Value = "listings.csv, CTA, Metra"
$ws.Cells.Item(16, 3).Value = "listings.csv, Chicago Data Portal"
$ws.Cells.Item(17, 3).Value = "listings.csv, Chicago Data Portal"
$ws.Cells.Item(19, 3).Value = "Crimes (Chicago Data Portal)"

# Widen column C to fit the longer source strings
# (ColumnWidth uses character units offset from the stored sheet "width" by
# the 5px/MDW padding Excel applies; 26.6667 round-trips to a stored 27.5.)
$ws.Columns.Item(3).ColumnWidth = 26.666666666666668

# Update selection / active cell to match the final saved state
$ws.Range("D26").Select()

# Page setup (paper size / orientation) picked up by the sheet on this save
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
